$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 19
$ws.Range("I8").Value = 19
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 57
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 82
$ws.Range("N8").ClearContents()

$ws.Range("H19").Value = 1103.25
$ws.Range("I19").Value = 990
$ws.Range("J19").Value = 1119.4286
$ws.Range("K19").Value = 990
$ws.Range("L19").Value = 1119.4286
$ws.Range("M19").Value = -815
$ws.Range("N19").Value = -1469.4286

$ws.Range("H70").Value = 2106.5334
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 1409.8
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 4229.4
$ws.Range("M70").Value = -10230
$ws.Range("N70").Value = -4769.4

$ws.Range("H73").Value = 2106.5334
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 1409.8
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 4229.4
$ws.Range("M73").Value = -9564
$ws.Range("N73").Value = -6101.4

$ws.Range("H126").Value = 84000
$ws.Range("J126").Value = 84000
$ws.Range("L126").Value = 84000
$ws.Range("N126").Value = -93880

$ws.Range("H132").Value = 901
$ws.Range("I132").Value = 759.7778
$ws.Range("J132").Value = 1324.6666
$ws.Range("K132").Value = 2279.3334
$ws.Range("L132").Value = 3973.9998
$ws.Range("M132").Value = 250.6666
$ws.Range("N132").Value = -9033.9998

$ws.Range("H138").Value = 2012.4166
$ws.Range("J138").Value = 2400
$ws.Range("L138").Value = 7200
$ws.Range("N138").Value = -17480

$ws.Range("H141").Value = 708.3333
$ws.Range("I141").Value = 734.375
$ws.Range("J141").Value = 500
$ws.Range("K141").Value = 2203.125
$ws.Range("L141").Value = 1500
$ws.Range("M141").Value = 2976.875
$ws.Range("N141").Value = -11860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 76747.25
$ws.Range("I124").Value = 74996
$ws.Range("K124").Value = 74996
$ws.Range("M124").Value = -70086

$ws.Range("H132").Value = 8614
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H75").Value = 5333
$ws.Range("J75").Value = 6000
$ws.Range("L75").Value = 6000
$ws.Range("N75").Value = -7872

$ws.Range("H78").Value = 5333
$ws.Range("J78").Value = 6000
$ws.Range("L78").Value = 18000
$ws.Range("N78").Value = -27360

$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1398.5
$ws.Range("I31").Value = 1398.5
$ws.Range("K31").Value = 1398.5
$ws.Range("M31").Value = -1103.5

$ws.Range("H34").Value = 1398.5
$ws.Range("I34").Value = 1398.5
$ws.Range("K34").Value = 1398.5
$ws.Range("M34").Value = -1196.5

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H99").Value = 1252499
$ws.Range("I99").Value = 3332
$ws.Range("K99").Value = 3332
$ws.Range("M99").Value = -1834

$ws.Range("H126").Value = 1252499
$ws.Range("I126").Value = 3332
$ws.Range("K126").Value = 9996
$ws.Range("M126").Value = -7526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 235.2
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws.Range("H47").Value = 669.6667
$ws.Range("I47").Value = 669.6667
$ws.Range("K47").Value = 2009.0001
$ws.Range("M47").Value = -1578.0001

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -35242

$ws.Range("H120").Value = 21015
$ws.Range("I120").Value = 2030
$ws.Range("J120").Value = 40000
$ws.Range("K120").Value = 6090
$ws.Range("L120").Value = 120000
$ws.Range("M120").Value = -1252
$ws.Range("N120").Value = -129676

$ws.Range("H134").Value = 4822
$ws.Range("I134").Value = 4822
$ws.Range("K134").Value = 14466
$ws.Range("M134").Value = -9396

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1067423.9
$ws.Range("I3").Value = 1000251.5
$ws.Range("J3").Value = 1112205.5
$ws.Range("K3").Value = 1000251.5
$ws.Range("L3").Value = 1112205.5
$ws.Range("M3").Value = -1000135.5
$ws.Range("N3").Value = -1112437.5

$ws.Range("H13").Value = 1133.3334
$ws.Range("I13").Value = 1133.3334
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1133.3334
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -994.3334
$ws.Range("N13").ClearContents()

$ws.Range("H111").Value = 79995
$ws.Range("J111").Value = 79995
$ws.Range("L111").Value = 79995
$ws.Range("N111").Value = -86129

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H93").Value = 943.75
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
